$d = $word.ActiveDocument

# 1. Fix the double space "pick up the  parrot" -> "pick up the parrot"
$find = $d.Content.Find
$find.Execute("pick up the  parrot", $true, $false, $false, $false, $false, `
               $true, 1, $false, "pick up the parrot", 2) | Out-Null

# 2. Insert three new paragraphs after the paragraph ending "...opposite side."
#    (the paragraph that also holds the _GoBack bookmark), each indented
#    with ind w:left=360, matching the Problem 1 section's trailing blank
#    line + new "Problem 2" heading + its first numbered item.
$pTarget = $d.Paragraphs.Item(8)
$rTarget = $pTarget.Range

# 2a. Blank spacer paragraph
$rTarget.InsertParagraphAfter() | Out-Null
$pBlank = $d.Paragraphs.Item(9)
$pBlank.Range.ListFormat.RemoveNumbers()
$pBlank.Style = "Normal"
$pBlank.LeftIndent = 18

# 2b. "Problem 2: Socks in the Dark" heading paragraph
$pBlank.Range.InsertParagraphAfter() | Out-Null
$pHeading = $d.Paragraphs.Item(10)
$pHeading.Range.ListFormat.RemoveNumbers()
$pHeading.Style = "Normal"
$pHeading.LeftIndent = 18
$pHeading.Range.Text = "Problem 2: Socks in the Dark"

# 2c. First body paragraph describing the problem
$pHeading.Range.InsertParagraphAfter() | Out-Null
$pBody = $d.Paragraphs.Item(11)
$pBody.Range.ListFormat.RemoveNumbers()
$pBody.Style = "Normal"
$pBody.LeftIndent = 18
$pBody.Range.Text = "1. So the problem at hand is that I need to pull socks from a drawer but can’t see the selection until I have pulled them because it is dark. There are 20 socks in total with 10 of them being black, 6 being brown, and 4 being white. The goal is for me to figure out the least amount of socks I would have to pull in order to get a) one matching pair and b) one matching pair of each color."

# 3. Add a lastRenderedPageBreak before the trailing tab run
$find2 = $d.Content.Find
$find2.Execute("^t", $true, $false, $false, $false, $false, `
                $true, 1, $false, "", 0) | Out-Null
